$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Clear existing data rows (below header) before rewriting, since row count changes and order changes
$ws1.Range("A2:C49").ClearContents()
$ws2.Range("A2:C6").ClearContents()

# Sheet 1 ("新增翻译") data rows
$ws1.Cells.Item(2, 1).Value = "ilylia"
$ws1.Cells.Item(2, 2).Value = 15015
$ws1.Cells.Item(2, 3).Value = 13
$ws1.Cells.Item(3, 1).Value = "huntsman-li"
$ws1.Cells.Item(3, 2).Value = 11209
$ws1.Cells.Item(3, 3).Value = 11
$ws1.Cells.Item(4, 1).Value = "JasonRD"
$ws1.Cells.Item(4, 2).Value = 5795
$ws1.Cells.Item(4, 3).Value = 5
$ws1.Cells.Item(5, 1).Value = "dk-lockdown"
$ws1.Cells.Item(5, 2).Value = 3434
$ws1.Cells.Item(5, 3).Value = 5
$ws1.Cells.Item(6, 1).Value = "yuer1727"
$ws1.Cells.Item(6, 2).Value = 3413
$ws1.Cells.Item(6, 3).Value = 3
$ws1.Cells.Item(7, 1).Value = "luxious"
$ws1.Cells.Item(7, 2).Value = 3278
$ws1.Cells.Item(7, 3).Value = 5
$ws1.Cells.Item(8, 1).Value = "yuxiaobo96"
$ws1.Cells.Item(8, 2).Value = 3182
$ws1.Cells.Item(8, 3).Value = 17
$ws1.Cells.Item(9, 1).Value = "vivian99-wu"
$ws1.Cells.Item(9, 2).Value = 3076
$ws1.Cells.Item(9, 3).Value = 3
$ws1.Cells.Item(10, 1).Value = "cocotyty"
$ws1.Cells.Item(10, 2).Value = 2765
$ws1.Cells.Item(10, 3).Value = 4
$ws1.Cells.Item(11, 1).Value = "ExBs2724"
$ws1.Cells.Item(11, 2).Value = 2636
$ws1.Cells.Item(11, 3).Value = 6
$ws1.Cells.Item(12, 1).Value = "tanjunchen"
$ws1.Cells.Item(12, 2).Value = 2520
$ws1.Cells.Item(12, 3).Value = 5
$ws1.Cells.Item(13, 1).Value = "idefav"
$ws1.Cells.Item(13, 2).Value = 2423
$ws1.Cells.Item(13, 3).Value = 2
$ws1.Cells.Item(14, 1).Value = "niceforbear"
$ws1.Cells.Item(14, 2).Value = 2282
$ws1.Cells.Item(14, 3).Value = 5
$ws1.Cells.Item(15, 1).Value = "gorda"
$ws1.Cells.Item(15, 2).Value = 2170
$ws1.Cells.Item(15, 3).Value = 3
$ws1.Cells.Item(16, 1).Value = "Lovnx"
$ws1.Cells.Item(16, 2).Value = 1995
$ws1.Cells.Item(16, 3).Value = 7
$ws1.Cells.Item(17, 1).Value = "GuangmingLuo"
$ws1.Cells.Item(17, 2).Value = 1861
$ws1.Cells.Item(17, 3).Value = 1
$ws1.Cells.Item(18, 1).Value = "zqzzq"
$ws1.Cells.Item(18, 2).Value = 1858
$ws1.Cells.Item(18, 3).Value = 5
$ws1.Cells.Item(19, 1).Value = "malphi"
$ws1.Cells.Item(19, 2).Value = 1606
$ws1.Cells.Item(19, 3).Value = 3
$ws1.Cells.Item(20, 1).Value = "sirius1024"
$ws1.Cells.Item(20, 2).Value = 1497
$ws1.Cells.Item(20, 3).Value = 2
$ws1.Cells.Item(21, 1).Value = "wangzewang"
$ws1.Cells.Item(21, 2).Value = 1412
$ws1.Cells.Item(21, 3).Value = 2
$ws1.Cells.Item(22, 1).Value = "rootsongjc"
$ws1.Cells.Item(22, 2).Value = 1272
$ws1.Cells.Item(22, 3).Value = 2
$ws1.Cells.Item(23, 1).Value = "WisWang"
$ws1.Cells.Item(23, 2).Value = 1185
$ws1.Cells.Item(23, 3).Value = 4
$ws1.Cells.Item(24, 1).Value = "arunfung"
$ws1.Cells.Item(24, 2).Value = 1095
$ws1.Cells.Item(24, 3).Value = 4
$ws1.Cells.Item(25, 1).Value = "JHDST"
$ws1.Cells.Item(25, 2).Value = 996
$ws1.Cells.Item(25, 3).Value = 6
$ws1.Cells.Item(26, 1).Value = "wenhuwang"
$ws1.Cells.Item(26, 2).Value = 974
$ws1.Cells.Item(26, 3).Value = 2
$ws1.Cells.Item(27, 1).Value = "homilly"
$ws1.Cells.Item(27, 2).Value = 880
$ws1.Cells.Item(27, 3).Value = 1
$ws1.Cells.Item(28, 1).Value = "tolbkni"
$ws1.Cells.Item(28, 2).Value = 877
$ws1.Cells.Item(28, 3).Value = 1
$ws1.Cells.Item(29, 1).Value = "jakeslee"
$ws1.Cells.Item(29, 2).Value = 831
$ws1.Cells.Item(29, 3).Value = 3
$ws1.Cells.Item(30, 1).Value = "gaohuag"
$ws1.Cells.Item(30, 2).Value = 726
$ws1.Cells.Item(30, 3).Value = 1
$ws1.Cells.Item(31, 1).Value = "zyt312074545"
$ws1.Cells.Item(31, 2).Value = 719
$ws1.Cells.Item(31, 3).Value = 3
$ws1.Cells.Item(32, 1).Value = "csdnshyang"
$ws1.Cells.Item(32, 2).Value = 628
$ws1.Cells.Item(32, 3).Value = 1
$ws1.Cells.Item(33, 1).Value = "yuxiaoba"
$ws1.Cells.Item(33, 2).Value = 513
$ws1.Cells.Item(33, 3).Value = 1
$ws1.Cells.Item(34, 1).Value = "zxh326"
$ws1.Cells.Item(34, 2).Value = 468
$ws1.Cells.Item(34, 3).Value = 5
$ws1.Cells.Item(35, 1).Value = "mrshengzyzy"
$ws1.Cells.Item(35, 2).Value = 456
$ws1.Cells.Item(35, 3).Value = 3
$ws1.Cells.Item(36, 1).Value = "JuwanXu"
$ws1.Cells.Item(36, 2).Value = 395
$ws1.Cells.Item(36, 3).Value = 12
$ws1.Cells.Item(37, 1).Value = "shicheng0829"
$ws1.Cells.Item(37, 2).Value = 381
$ws1.Cells.Item(37, 3).Value = 1
$ws1.Cells.Item(38, 1).Value = "zzzhy"
$ws1.Cells.Item(38, 2).Value = 362
$ws1.Cells.Item(38, 3).Value = 1
$ws1.Cells.Item(39, 1).Value = "AsCat"
$ws1.Cells.Item(39, 2).Value = 228
$ws1.Cells.Item(39, 3).Value = 1
$ws1.Cells.Item(40, 1).Value = "TomatoAres"
$ws1.Cells.Item(40, 2).Value = 197
$ws1.Cells.Item(40, 3).Value = 2
$ws1.Cells.Item(41, 1).Value = "lilinji"
$ws1.Cells.Item(41, 2).Value = 189
$ws1.Cells.Item(41, 3).Value = 2
$ws1.Cells.Item(42, 1).Value = "qunqiang"
$ws1.Cells.Item(42, 2).Value = 179
$ws1.Cells.Item(42, 3).Value = 3
$ws1.Cells.Item(43, 1).Value = "sunny0826"
$ws1.Cells.Item(43, 2).Value = 163
$ws1.Cells.Item(43, 3).Value = 1
$ws1.Cells.Item(44, 1).Value = "lengrongfu"
$ws1.Cells.Item(44, 2).Value = 147
$ws1.Cells.Item(44, 3).Value = 2
$ws1.Cells.Item(45, 1).Value = "wuti1609"
$ws1.Cells.Item(45, 2).Value = 144
$ws1.Cells.Item(45, 3).Value = 1
$ws1.Cells.Item(46, 1).Value = "dotw"
$ws1.Cells.Item(46, 2).Value = 126
$ws1.Cells.Item(46, 3).Value = 1
$ws1.Cells.Item(47, 1).Value = "kylesliu"
$ws1.Cells.Item(47, 2).Value = 96
$ws1.Cells.Item(47, 3).Value = 2
$ws1.Cells.Item(48, 1).Value = "GanymedeNil"
$ws1.Cells.Item(48, 2).Value = 79
$ws1.Cells.Item(48, 3).Value = 1
$ws1.Cells.Item(49, 1).Value = "innerpeacez"
$ws1.Cells.Item(49, 2).Value = 76
$ws1.Cells.Item(49, 3).Value = 1
$ws1.Cells.Item(50, 1).Value = "Luluda"
$ws1.Cells.Item(50, 2).Value = 41
$ws1.Cells.Item(50, 3).Value = 1
$ws1.Cells.Item(51, 1).Value = "vflong"
$ws1.Cells.Item(51, 2).Value = 37
$ws1.Cells.Item(51, 3).Value = 2
$ws1.Cells.Item(52, 1).Value = "5idu"
$ws1.Cells.Item(52, 2).Value = 32
$ws1.Cells.Item(52, 3).Value = 1
$ws1.Cells.Item(53, 1).Value = "koonchen"
$ws1.Cells.Item(53, 2).Value = 22
$ws1.Cells.Item(53, 3).Value = 1
$ws1.Cells.Item(54, 1).Value = "sniperking1234"
$ws1.Cells.Item(54, 2).Value = 10
$ws1.Cells.Item(54, 3).Value = 1

# Sheet 2 ("更新翻译") data rows
$ws2.Cells.Item(2, 1).Value = "tanjunchen"
$ws2.Cells.Item(2, 2).Value = 13506
$ws2.Cells.Item(2, 3).Value = 5
$ws2.Cells.Item(3, 1).Value = "hwdef"
$ws2.Cells.Item(3, 2).Value = 1538
$ws2.Cells.Item(3, 3).Value = 1
$ws2.Cells.Item(4, 1).Value = "ExBs2724"
$ws2.Cells.Item(4, 2).Value = 588
$ws2.Cells.Item(4, 3).Value = 6
$ws2.Cells.Item(5, 1).Value = "vivian99-wu"
$ws2.Cells.Item(5, 2).Value = 506
$ws2.Cells.Item(5, 3).Value = 3
$ws2.Cells.Item(6, 1).Value = "zyt312074545"
$ws2.Cells.Item(6, 2).Value = 375
$ws2.Cells.Item(6, 3).Value = 3
$ws2.Cells.Item(7, 1).Value = "JasonRD"
$ws2.Cells.Item(7, 2).Value = 257
$ws2.Cells.Item(7, 3).Value = 5
$ws2.Cells.Item(8, 1).Value = "ilylia"
$ws2.Cells.Item(8, 2).Value = 222
$ws2.Cells.Item(8, 3).Value = 13
$ws2.Cells.Item(9, 1).Value = "Lovnx"
$ws2.Cells.Item(9, 2).Value = 211
$ws2.Cells.Item(9, 3).Value = 7
$ws2.Cells.Item(10, 1).Value = "idefav"
$ws2.Cells.Item(10, 2).Value = 208
$ws2.Cells.Item(10, 3).Value = 2
$ws2.Cells.Item(11, 1).Value = "WisWang"
$ws2.Cells.Item(11, 2).Value = 206
$ws2.Cells.Item(11, 3).Value = 4
$ws2.Cells.Item(12, 1).Value = "gorda"
$ws2.Cells.Item(12, 2).Value = 195
$ws2.Cells.Item(12, 3).Value = 3
$ws2.Cells.Item(13, 1).Value = "dk-lockdown"
$ws2.Cells.Item(13, 2).Value = 125
$ws2.Cells.Item(13, 3).Value = 5
$ws2.Cells.Item(14, 1).Value = "rootsongjc"
$ws2.Cells.Item(14, 2).Value = 24
$ws2.Cells.Item(14, 3).Value = 2

Write-Output "done"
